$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-28 Monday", "2025-07-29 Tuesday"),
    @("72÷5=14, 2", "59÷8=7, 3"),
    @("73÷3=24, 1", "45÷7=6, 3"),
    @("68÷7=9, 5", "34÷8=4, 2"),
    @("81÷8=10, 1", "38÷2=19, 0"),
    @("93÷2=46, 1", "50÷3=16, 2"),
    @("25÷8=3, 1", "45÷5=9, 0"),
    @("23÷2=11, 1", "47÷2=23, 1"),
    @("98÷6=16, 2", "56÷3=18, 2"),
    @("99÷3=33, 0", "10÷7=1, 3"),
    @("25÷2=12, 1", "32÷6=5, 2"),
    @("58÷6=9, 4", "26÷6=4, 2"),
    @("29÷9=3, 2", "29÷5=5, 4"),
    @("92÷6=15, 2", "44÷9=4, 8"),
    @("83÷7=11, 6", "66÷2=33, 0"),
    @("51÷8=6, 3", "14÷2=7, 0"),
    @("42÷5=8, 2", "35÷6=5, 5"),
    @("57÷8=7, 1", "79÷7=11, 2"),
    @("56÷2=28, 0", "95÷4=23, 3"),
    @("31÷3=10, 1", "77÷7=11, 0"),
    @("98÷5=19, 3", "38÷7=5, 3"),
    @("66÷3=22, 0", "41÷8=5, 1"),
    @("16÷3=5, 1", "90÷8=11, 2"),
    @("80÷2=40, 0", "90÷4=22, 2"),
    @("18÷9=2, 0", "43÷7=6, 1"),
    @("58÷5=11, 3", "23÷2=11, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
